$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the match data that currently lives in rows 47 and 48.
#    Columns A-D, G, K, O, S stay put (they are identical between the two
#    rows); columns F, H, I, J, L, M, N, P, Q, R, T, U, V need to be
#    exchanged between row 47 and row 48.
# ---------------------------------------------------------------------------
$cols = @("F","H","I","J","L","M","N","P","Q","R","T","U","V")

foreach ($col in $cols) {
    $r47 = $col + "47"
    $r48 = $col + "48"
    $v47 = $ws.Range($r47).Value()
    $v48 = $ws.Range($r48).Value()
    $ws.Range($r47).Value = $v48
    $ws.Range($r48).Value = $v47
}

# ---------------------------------------------------------------------------
# 2) Append a new match (row 62), matching the formatting used by the
#    existing data rows (row 61 is used as the formatting template).
# ---------------------------------------------------------------------------
$ws.Range("A61:V61").Copy()
$ws.Range("A62").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = "morocco"
$ws.Cells.Item(62, 3).Value = "botola-pro"
$ws.Cells.Item(62, 4).Value = "2023-2024"
$ws.Cells.Item(62, 5).Value = 45236.70833333334
$ws.Cells.Item(62, 6).Value = "Raja Casablanca"
$ws.Cells.Item(62, 7).Value = 2
$ws.Cells.Item(62, 8).Value = "FUS Rabat"
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 1.96
$ws.Cells.Item(62, 11).Value = "05/11/2023 05:12"
$ws.Cells.Item(62, 12).Value = 2.02
$ws.Cells.Item(62, 13).Value = "06/11/2023 16:59"
$ws.Cells.Item(62, 14).Value = 2.94
$ws.Cells.Item(62, 15).Value = "05/11/2023 05:12"
$ws.Cells.Item(62, 16).Value = 2.86
$ws.Cells.Item(62, 17).Value = "06/11/2023 16:50"
$ws.Cells.Item(62, 18).Value = 3.78
$ws.Cells.Item(62, 19).Value = "05/11/2023 05:12"
$ws.Cells.Item(62, 20).Value = 4.42
$ws.Cells.Item(62, 21).Value = "06/11/2023 16:59"
$ws.Cells.Item(62, 22).Value = "https://www.betexplorer.com/football/morocco/botola-pro/raja-casablanca-fus-rabat/6y6heMDc/"
